# Update gh-pages output data (commit 456a3b4)
# Applies updated "想去人数" (F) / "最低票价" (G) values across the
# 展览, 演出 and 全部类型 worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3005
$ws.Range("G2").Value = 80
$ws.Range("F4").Value = 2551
$ws.Range("F7").Value = 86
$ws.Range("F9").Value = 3145
$ws.Range("F10").Value = 368
$ws.Range("F12").Value = 7687
$ws.Range("F13").Value = 382
$ws.Range("F20").Value = 9446
$ws.Range("F27").Value = 129
$ws.Range("F28").Value = 131
$ws.Range("F29").Value = 35
$ws.Range("F33").Value = 2624
$ws.Range("F36").Value = 2052
$ws.Range("F39").Value = 3968
$ws.Range("F40").Value = 218
$ws.Range("F41").Value = 48
$ws.Range("F43").Value = 113
$ws.Range("F44").Value = 256
$ws.Range("F45").Value = 65
$ws.Range("F47").Value = 72
$ws.Range("F49").Value = 67

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 159
$ws.Range("F16").Value = 13
$ws.Range("F21").Value = 14

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3005
$ws.Range("G3").Value = 80
$ws.Range("F7").Value = 2551
$ws.Range("F8").Value = 159
$ws.Range("F11").Value = 86
$ws.Range("F13").Value = 3145
$ws.Range("F14").Value = 368
$ws.Range("F18").Value = 7687
$ws.Range("F19").Value = 382
$ws.Range("F24").Value = 9446
$ws.Range("F29").Value = 129
$ws.Range("F30").Value = 131
$ws.Range("F31").Value = 35
$ws.Range("F35").Value = 2624
$ws.Range("F36").Value = 2052
$ws.Range("F40").Value = 3968
$ws.Range("F41").Value = 218
$ws.Range("F42").Value = 48
$ws.Range("F44").Value = 113
$ws.Range("F45").Value = 256
$ws.Range("F46").Value = 65
$ws.Range("F47").Value = 72
$ws.Range("F49").Value = 67
